$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing 3 rows (former ECs-sender block); new dimension becomes A1:T7
$ws.Range("A8:A10").EntireRow.Delete()

# Rewrite rows 2-7 with the refreshed TPM-derived values

# Row 2: FAPs -> ECs (Angpt1-Itgb1)
$ws.Cells.Item(2,1).Value2 = "FAPs"
$ws.Cells.Item(2,2).Value2 = "Angpt1"
$ws.Cells.Item(2,3).Value2 = "Itgb1"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 11.05178533333333
$ws.Cells.Item(2,8).Value2 = 33.155356
$ws.Cells.Item(2,9).Value2 = 0.9017494976312432
$ws.Cells.Item(2,10).Value2 = 0.9017494976312432
$ws.Cells.Item(2,11).Value2 = 3
$ws.Cells.Item(2,12).Value2 = 1
$ws.Cells.Item(2,13).Value2 = 61.04160633333334
$ws.Cells.Item(2,14).Value2 = 183.124819
$ws.Cells.Item(2,15).Value2 = 0.2043613460574534
$ws.Cells.Item(2,16).Value2 = 0.2043613460574534
$ws.Cells.Item(2,17).Value2 = 674.6187295978406
$ws.Cells.Item(2,18).Value2 = 6071.568566380564
$ws.Cells.Item(2,19).Value2 = 0.1842827411425533
$ws.Cells.Item(2,20).Value2 = 0.1842827411425533

# Row 3: FAPs -> FAPs (Angpt1-Itgb1)
$ws.Cells.Item(3,1).Value2 = "FAPs"
$ws.Cells.Item(3,2).Value2 = "Angpt1"
$ws.Cells.Item(3,3).Value2 = "Itgb1"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 11.05178533333333
$ws.Cells.Item(3,8).Value2 = 33.155356
$ws.Cells.Item(3,9).Value2 = 0.9017494976312432
$ws.Cells.Item(3,10).Value2 = 0.9017494976312432
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 106.314466
$ws.Cells.Item(3,14).Value2 = 318.943398
$ws.Cells.Item(3,15).Value2 = 0.3559304658284363
$ws.Cells.Item(3,16).Value2 = 0.3559304658284363
$ws.Cells.Item(3,17).Value2 = 1174.964656059965
$ws.Cells.Item(3,18).Value2 = 10574.68190453969
$ws.Cells.Item(3,19).Value2 = 0.3209601187524468
$ws.Cells.Item(3,20).Value2 = 0.3209601187524468

# Row 4: FAPs -> MuSCs (Angpt1-Itgb1)
$ws.Cells.Item(4,1).Value2 = "FAPs"
$ws.Cells.Item(4,2).Value2 = "Angpt1"
$ws.Cells.Item(4,3).Value2 = "Itgb1"
$ws.Cells.Item(4,4).Value2 = "MuSCs"
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,6).Value2 = 1
$ws.Cells.Item(4,7).Value2 = 11.05178533333333
$ws.Cells.Item(4,8).Value2 = 33.155356
$ws.Cells.Item(4,9).Value2 = 0.9017494976312432
$ws.Cells.Item(4,10).Value2 = 0.9017494976312432
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 131.3384093333333
$ws.Cells.Item(4,14).Value2 = 394.015228
$ws.Cells.Item(4,15).Value2 = 0.4397081881141102
$ws.Cells.Item(4,16).Value2 = 0.4397081881141103
$ws.Cells.Item(4,17).Value2 = 1451.523905973463
$ws.Cells.Item(4,18).Value2 = 13063.71515376117
$ws.Cells.Item(4,19).Value2 = 0.3965066377362431
$ws.Cells.Item(4,20).Value2 = 0.3965066377362431

# Row 5: MuSCs -> ECs (Angpt1-Itgb1)
$ws.Cells.Item(5,1).Value2 = "MuSCs"
$ws.Cells.Item(5,2).Value2 = "Angpt1"
$ws.Cells.Item(5,3).Value2 = "Itgb1"
$ws.Cells.Item(5,4).Value2 = "ECs"
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,6).Value2 = 1
$ws.Cells.Item(5,7).Value2 = 1.204152
$ws.Cells.Item(5,8).Value2 = 3.612456
$ws.Cells.Item(5,9).Value2 = 0.09825050236875665
$ws.Cells.Item(5,10).Value2 = 0.09825050236875667
$ws.Cells.Item(5,11).Value2 = 3
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 61.04160633333334
$ws.Cells.Item(5,14).Value2 = 183.124819
$ws.Cells.Item(5,15).Value2 = 0.2043613460574534
$ws.Cells.Item(5,16).Value2 = 0.2043613460574534
$ws.Cells.Item(5,17).Value2 = 73.503372349496
$ws.Cells.Item(5,18).Value2 = 661.5303511454641
$ws.Cells.Item(5,19).Value2 = 0.02007860491490013
$ws.Cells.Item(5,20).Value2 = 0.02007860491490013

# Row 6: MuSCs -> FAPs (Angpt1-Itgb1)
$ws.Cells.Item(6,1).Value2 = "MuSCs"
$ws.Cells.Item(6,2).Value2 = "Angpt1"
$ws.Cells.Item(6,3).Value2 = "Itgb1"
$ws.Cells.Item(6,4).Value2 = "FAPs"
$ws.Cells.Item(6,5).Value2 = 3
$ws.Cells.Item(6,6).Value2 = 1
$ws.Cells.Item(6,7).Value2 = 1.204152
$ws.Cells.Item(6,8).Value2 = 3.612456
$ws.Cells.Item(6,9).Value2 = 0.09825050236875665
$ws.Cells.Item(6,10).Value2 = 0.09825050236875667
$ws.Cells.Item(6,11).Value2 = 3
$ws.Cells.Item(6,12).Value2 = 1
$ws.Cells.Item(6,13).Value2 = 106.314466
$ws.Cells.Item(6,14).Value2 = 318.943398
$ws.Cells.Item(6,15).Value2 = 0.3559304658284363
$ws.Cells.Item(6,16).Value2 = 0.3559304658284363
$ws.Cells.Item(6,17).Value2 = 128.018776862832
$ws.Cells.Item(6,18).Value2 = 1152.168991765488
$ws.Cells.Item(6,19).Value2 = 0.03497034707598944
$ws.Cells.Item(6,20).Value2 = 0.03497034707598945

# Row 7: MuSCs -> MuSCs (Angpt1-Itgb1)
$ws.Cells.Item(7,1).Value2 = "MuSCs"
$ws.Cells.Item(7,2).Value2 = "Angpt1"
$ws.Cells.Item(7,3).Value2 = "Itgb1"
$ws.Cells.Item(7,4).Value2 = "MuSCs"
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,6).Value2 = 1
$ws.Cells.Item(7,7).Value2 = 1.204152
$ws.Cells.Item(7,8).Value2 = 3.612456
$ws.Cells.Item(7,9).Value2 = 0.09825050236875665
$ws.Cells.Item(7,10).Value2 = 0.09825050236875667
$ws.Cells.Item(7,11).Value2 = 3
$ws.Cells.Item(7,12).Value2 = 1
$ws.Cells.Item(7,13).Value2 = 131.3384093333333
$ws.Cells.Item(7,14).Value2 = 394.015228
$ws.Cells.Item(7,15).Value2 = 0.4397081881141102
$ws.Cells.Item(7,16).Value2 = 0.4397081881141103
$ws.Cells.Item(7,17).Value2 = 158.151408275552
$ws.Cells.Item(7,18).Value2 = 1423.362674479968
$ws.Cells.Item(7,19).Value2 = 0.04320155037786708
$ws.Cells.Item(7,20).Value2 = 0.0432015503778671
